# epexspot_prices.xlsx automated daily update
# - "Prix Spot": append a new date column (AJ) = "19-jul" with one price per hour row
# - "Gaz": append a new date row (33) = 2025-07-17
# - "CO2":  append a new date row (33) = 2025-07-17

$wb = $excel.ActiveWorkbook

# xlPasteValues / xlPasteFormats constants used below
$xlPasteValues  = -4163
$xlPasteFormats = -4122

# ---------------------------------------------------------------
# Sheet "Prix Spot": new column AJ ("19-jul")
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Prix Spot")

# Seed AJ1 from AI1 so it inherits the bold/bordered header style,
# then swap the text in place (formula -> paste-as-value keeps the
# style untouched and avoids minting a brand new cell format).
$ws.Range("AI1").Copy($ws.Range("AJ1"))
$ws.Range("AJ1").Formula = '="19-jul"'
$ws.Range("AJ1").Copy()
$ws.Range("AJ1").PasteSpecial($xlPasteValues)

$ws.Range("AJ2").Value  = 98.8
$ws.Range("AJ3").Value  = 90.79000000000001
$ws.Range("AJ4").Value  = 79.70999999999999
$ws.Range("AJ5").Value  = 53.41
$ws.Range("AJ6").Value  = 49.57
$ws.Range("AJ7").Value  = 48.62
$ws.Range("AJ8").Value  = 47.99
$ws.Range("AJ9").Value  = 56.25
$ws.Range("AJ10").Value = 56.62
$ws.Range("AJ11").Value = 34.46
$ws.Range("AJ12").Value = 10
$ws.Range("AJ13").Value = 0.65
$ws.Range("AJ14").Value = 0
$ws.Range("AJ15").Value = 0
$ws.Range("AJ16").Value = 0
$ws.Range("AJ17").Value = 2.37
$ws.Range("AJ18").Value = 18.82
$ws.Range("AJ19").Value = 19.38
$ws.Range("AJ20").Value = 30.16
$ws.Range("AJ21").Value = 42.32
$ws.Range("AJ22").Value = 72.92
$ws.Range("AJ23").Value = 94.17
$ws.Range("AJ24").Value = 106.8
$ws.Range("AJ25").Value = 96.25

$excel.CutCopyMode = $false

# ---------------------------------------------------------------
# Sheet "Gaz": new row 33 (2025-07-17 / 33.6)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Gaz")

# Route the date string through a formula + paste-as-value so Excel's
# "looks like a date" auto-conversion never kicks in, while leaving
# the cell on the default (unstyled) format like its neighbours.
$ws.Range("A33").Formula = '="2025-07-17"'
$ws.Range("A33").Copy()
$ws.Range("A33").PasteSpecial($xlPasteValues)
$ws.Range("B33").Value = 33.6

$excel.CutCopyMode = $false

# ---------------------------------------------------------------
# Sheet "CO2": new row 33 (2025-07-17 / 69.81999999999999)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("CO2")

$ws.Range("A33").Formula = '="2025-07-17"'
$ws.Range("A33").Copy()
$ws.Range("A33").PasteSpecial($xlPasteValues)
$ws.Range("B33").Value = 69.81999999999999

$excel.CutCopyMode = $false
